$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns D, J, K, L, M, P for rows 2..70
# Each entry: row, D, J, K, L, M, P
$data = @(
    @(2, 44460, 800, 4000, 4500, 4250, 4250),
    @(3, 44335, 600, 3000, 3500, 3250, 3250),
    @(4, 44441, 600, 4500, 5000, 4750, 4750),
    @(5, 44504, 760, 3000, 3500, 3250, 3250),
    @(6, 44455, 600, 4500, 5000, 4750, 4750),
    @(7, 44334, 760, 3000, 3500, 3250, 3250),
    @(8, 44522, 800, 3000, 3500, 3250, 3250),
    @(9, 44453, 800, 4500, 5000, 4750, 4750),
    @(10, 44314, 800, 2500, 3000, 2750, 2750),
    @(11, 44462, 660, 4000, 4500, 4250, 4250),
    @(12, 44377, 600, 4000, 4500, 4250, 4250),
    @(13, 44497, 740, 3800, 4000, 3900, 3900),
    @(14, 44435, 1500, 4500, 5000, 4750, 4750),
    @(15, 44445, 600, 4500, 5000, 4750, 4750),
    @(16, 44356, 600, 3000, 3500, 3250, 3250),
    @(17, 44342, 560, 3000, 3500, 3250, 3250),
    @(18, 44509, 800, 3500, 4000, 3750, 3750),
    @(19, 44448, 640, 4500, 5000, 4750, 4750),
    @(20, 44515, 800, 3000, 4000, 3500, 3500),
    @(21, 44397, 800, 4000, 4500, 4250, 4250),
    @(22, 44523, 760, 3000, 4000, 3500, 3500),
    @(23, 44169, 2400, 3000, 3500, 3250, 3250),
    @(24, 44484, 840, 3500, 4000, 3750, 3750),
    @(25, 44530, 800, 3000, 4000, 3500, 3500),
    @(26, 44512, 800, 3000, 3500, 3250, 3250),
    @(27, 44537, 760, 3500, 4000, 3750, 3750),
    @(28, 44407, 720, 4000, 4500, 4250, 4250),
    @(29, 44315, 700, 2500, 3000, 2750, 2750),
    @(30, 44483, 700, 3500, 4000, 3750, 3750),
    @(31, 44505, 800, 3500, 4000, 3750, 3750),
    @(32, 44348, 700, 3000, 3500, 3250, 3250),
    @(33, 44488, 800, 3500, 4000, 3750, 3750),
    @(34, 44172, 760, 3000, 3500, 3250, 3250),
    @(35, 44162, 2000, 2800, 3000, 2900, 2900),
    @(36, 44349, 560, 3000, 3500, 3250, 3250),
    @(37, 44469, 700, 4000, 4500, 4250, 4250),
    @(38, 44525, 720, 3000, 4000, 3500, 3500),
    @(39, 44165, 1000, 3000, 3500, 3250, 3250),
    @(40, 44516, 740, 3000, 4000, 3500, 3500),
    @(41, 44427, 600, 4500, 5000, 4750, 4750),
    @(42, 44533, 900, 3000, 4000, 3500, 3500),
    @(43, 44176, 2000, 3000, 3500, 3250, 3250),
    @(44, 44539, 600, 3000, 4000, 3500, 3500),
    @(45, 44379, 800, 4000, 4500, 4250, 4250),
    @(46, 44532, 740, 3000, 4000, 3500, 3500),
    @(47, 44498, 900, 3800, 4000, 3900, 3900),
    @(48, 44449, 700, 4000, 4500, 4250, 4250),
    @(49, 44526, 800, 3000, 4000, 3500, 3500),
    @(50, 44425, 900, 4500, 5000, 4750, 4750),
    @(51, 44476, 600, 3500, 4000, 3750, 3750),
    @(52, 44418, 800, 4500, 5000, 4750, 4750),
    @(53, 44434, 600, 4500, 5000, 4750, 4750),
    @(54, 44467, 840, 4000, 4500, 4250, 4250),
    @(55, 44490, 660, 3500, 4000, 3750, 3750),
    @(56, 44341, 700, 3000, 3500, 3250, 3250),
    @(57, 44519, 800, 3500, 4000, 3750, 3750),
    @(58, 44420, 900, 4500, 5000, 4750, 4750),
    @(59, 44474, 760, 3500, 4000, 3750, 3750),
    @(60, 44446, 800, 4500, 5000, 4750, 4750),
    @(61, 44411, 880, 4000, 4500, 4250, 4250),
    @(62, 44432, 900, 4500, 5000, 4750, 4750),
    @(63, 44452, 600, 4500, 5000, 4750, 4750),
    @(64, 44473, 600, 3500, 4000, 3750, 3750),
    @(65, 44536, 1000, 3500, 4000, 3750, 3750),
    @(66, 44511, 760, 3000, 3500, 3250, 3250),
    @(67, 44463, 800, 4000, 4500, 4250, 4250),
    @(68, 44365, 800, 3500, 4000, 3750, 3750),
    @(69, 44518, 760, 3000, 4000, 3500, 3500),
    @(70, 44540, 500, 3000, 4000, 3500, 3500)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 4).Value = $row[1]
    $ws.Cells.Item($r, 10).Value = $row[2]
    $ws.Cells.Item($r, 11).Value = $row[3]
    $ws.Cells.Item($r, 12).Value = $row[4]
    $ws.Cells.Item($r, 13).Value = $row[5]
    $ws.Cells.Item($r, 16).Value = $row[6]
}

Write-Host "Applied updates to $($data.Count) rows"